$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "'8"
$ws.Range("B4").Value = "[BUG] deployment failing"
$ws.Range("C4").Value = "open"
$ws.Range("D4").Value = "2025-03-24T08:21:23Z"
$ws.Range("E4").Value = "bug"
